# Auto-update gym prices
$wb = $excel.ActiveWorkbook

# --- Sheet: "4x4 Squat Racks" ---
$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")

# C2: $2,142.00 -> $2,139.00  (force text so the leading "$" isn't
# reinterpreted as a currency number; reset style afterwards so no
# NumberFormat/style gets stamped onto the cell)
$ws1.Range("C2").NumberFormat = "@"
$ws1.Range("C2").Value = "$2,139.00"
$ws1.Range("C2").Style = "Normal"

# C3: $1,299.99 -> Price not available
$ws1.Range("C3").Value = "Price not available"

# --- Sheet: "Squat Stands" ---
$ws2 = $wb.Worksheets.Item("Squat Stands")

# C2: $1,549.00 -> $1,546.00
$ws2.Range("C2").NumberFormat = "@"
$ws2.Range("C2").Value = "$1,546.00"
$ws2.Range("C2").Style = "Normal"

# C3 (previously empty/missing): -> Price not available
$ws2.Range("C3").Value = "Price not available"

# --- Sheet: "Leg Extensions" ---
$ws3 = $wb.Worksheets.Item("Leg Extensions")

# A3: Leg Extension and Curl Machine | 10 - 250 LB Selector -> Unknown Product
$ws3.Range("A3").Value = "Unknown Product"

# C3: $2,909.99 -> Price not available
$ws3.Range("C3").Value = "Price not available"
